$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "BestBefore" in column D (becomes shared string index 6)
$ws.Range("D1").Value = "BestBefore"

# Add the date values (stored as raw date serials) for each data row
$ws.Range("D2").Value = 45615
$ws.Range("D3").Value = 45631
$ws.Range("D4").Value = 45665

# Apply a date number format (maps to builtin numFmtId 14) to D2,
# then copy that formatting down to D3:D4 so all three cells share a
# single cell-style (xf) entry instead of each getting its own.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3:D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Give column D an explicit width similar to the authored workbook
$ws.Columns.Item(4).ColumnWidth = 13.14

# Move/represent the active selection on the newly added cell, matching
# where the author was working when the change was made
$ws.Range("D2").Select()
